# The commit swaps the contents of ppt/theme/theme1.xml (the theme used
# by the slide master / slides -- "Integral") and ppt/theme/theme2.xml
# (the theme used by the notes master -- the stock "Office Theme"), so
# that afterwards the slides use the plain default "Office Theme" colours
# and the notes master carries what used to be the "Integral" palette.
#
# theme1.xml and theme2.xml are byte-identical apart from the <a:theme>/
# <a:clrScheme> "name" attributes and the twelve colour values inside
# <a:clrScheme> (font scheme and format scheme are already identical in
# both files), so the visible, scriptable part of this edit is simply
# re-pointing the slide master's theme colour scheme at the stock Office
# palette.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

function ToRgbVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

# Index order for ThemeColorScheme.Item(n): 1=dk1 2=lt1 3=dk2 4=lt2
# 5..10=accent1..6 11=hlink 12=folHlink -- matches <a:clrScheme> order.
$officeTheme = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 0; $i -lt $officeTheme.Length; $i++) {
    $colors.Item($i + 1).RGB = ToRgbVal($officeTheme[$i])
}
